# Updated cryptos list with GitHub Actions
# Applies the latest price/volume (and two re-ranked rows) scraped from
# coinranking.com to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while always keeping it as TEXT.
# The sheet stores Price/Volume figures as strings (e.g. "5.750", "28.607.47",
# "  +1.56%  ") even though many of them look numeric. A plain
# $ws.Range(...).Value = "5.750" assignment would make Excel's COM layer
# auto-convert such look-alike numbers into real numeric values (and e.g.
# drop the significant trailing zero), so numeric-looking strings are
# written with a leading apostrophe to force a text entry, exactly like
# typing '5.750 into the cell in the Excel UI.
function Set-TextValue {
    param($Range, $Value)
    if ($Value -match '^[0-9]+(\.[0-9]+)?$') {
        $Range.Value = "'" + $Value
    } else {
        $Range.Value = $Value
    }
}

# row -> column/value pairs taken from the diff
$rowUpdates = @(
    @{ Row = 2;  D = "28.607.47";   E = "  +1.56%  " },
    @{ Row = 3;  D = "1.827.44";    E = "  +1.25%  " },
    @{ Row = 4;  D = "1.002";       E = "  +0.07%  " },
    @{ Row = 5;  D = "316.73";      E = "  -0.01%  " },
    @{ Row = 7;  D = "0.5329";      E = "  -1.98%  " },
    @{ Row = 8;  D = "0.3972";      E = "  +4.78%  " },
    @{ Row = 9;  D = "0.07766";     E = "  +3.79%  " },
    @{ Row = 10; B = "OKB";     C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb";        D = "42.05";   E = "  -0.24%  " },
    @{ Row = 11; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic";   D = "1.119";   E = "  +2.03%  " },
    @{ Row = 12; D = "21.14";       E = "  +2.73%  " },
    @{ Row = 13; D = "6.322";       E = "  +1.96%  " },
    @{ Row = 14; D = "7.573";       E = "  +2.82%  " },
    @{ Row = 15; D = "1.002";       E = "  +0.10%  " },
    @{ Row = 16; D = "1.819.34";    E = "  +0.82%  " },
    @{ Row = 17; D = "93.53";       E = "  +3.83%  " },
    @{ Row = 18; D = "0.00001091";  E = "  +2.30%  " },
    @{ Row = 19; D = "0.06619";     E = "  +1.64%  " },
    @{ Row = 20; D = "17.80";       E = "  +1.90%  " },
    @{ Row = 21; E = "  +0.08%  " },
    @{ Row = 22; D = "6.083";       E = "  +2.47%  " },
    @{ Row = 23; D = "28.607.90";   E = "  +1.47%  " },
    @{ Row = 24; D = "11.20";       E = "  -0.17%  " },
    @{ Row = 25; D = "2.235";       E = "  +6.81%  " },
    @{ Row = 26; D = "20.81";       E = "  +1.43%  " },
    @{ Row = 27; D = "156.88";      E = "  +0.44%  " },
    @{ Row = 28; D = "2.423";       E = "  +3.19%  " },
    @{ Row = 29; D = "2.007.71";    E = "  -0.25%  " },
    @{ Row = 30; D = "125.24";      E = "  +2.49%  " },
    @{ Row = 31; D = "1.152";       E = "  +2.53%  " },
    @{ Row = 32; D = "0.1129";      E = "  +0.57%  " },
    @{ Row = 33; D = "5.750";       E = "  +2.77%  " },
    @{ Row = 34; D = "3.663";       E = "  +0.00%  " },
    @{ Row = 35; D = "0.07338";     E = "  +5.47%  " },
    @{ Row = 36; D = "0.2273";      E = "  +1.80%  " },
    @{ Row = 37; D = "0.02352";     E = "  +1.92%  " },
    @{ Row = 38; D = "8.918";       E = "  +5.06%  " },
    @{ Row = 39; D = "5.211";       E = "  +2.15%  " },
    @{ Row = 40; D = "11.43";       E = "  +2.21%  " },
    @{ Row = 41; D = "0.6307";      E = "  +1.89%  " },
    @{ Row = 42; D = "1.197";       E = "  +1.80%  " },
    @{ Row = 43; E = "  +0.04%  " },
    @{ Row = 44; D = "1.399";       E = "  -1.47%  " },
    @{ Row = 45; D = "13.63";       E = "  +1.73%  " },
    @{ Row = 46; D = "0.5950" },
    @{ Row = 47; D = "3.722";       E = "  +0.95%  " },
    @{ Row = 48; D = "125.59";      E = "  +0.46%  " },
    @{ Row = 49; D = "2.000";       E = "  +3.71%  " },
    @{ Row = 50; E = "  +0.30%  " },
    @{ Row = 51; D = "0.06965";     E = "  +2.09%  " }
)

foreach ($update in $rowUpdates) {
    $row = $update.Row
    if ($update.ContainsKey("B")) {
        $ws.Range("B$row").Value = $update.B
    }
    if ($update.ContainsKey("C")) {
        $ws.Range("C$row").Value = $update.C
    }
    if ($update.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$row") $update.D
    }
    if ($update.ContainsKey("E")) {
        $ws.Range("E$row").Value = $update.E
    }
}
